# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.901.27'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '3.857.42'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'598.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").Value = "'166.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.21%  '
$ws.Range("D7").Value = '3.858.30'
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").Value = "'6.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = "'0.456"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = "'36.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").Value = '4.503.43'
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").Value = '3.878.68'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '67.949.58'
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("E18").Value = '  +6.93%  '
$ws.Range("D19").Value = "'7.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("E20").Value = '  -1.39%  '
$ws.Range("D21").Value = "'10.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").Value = "'464.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.79%  '
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("D24").Value = "'0.0000163"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").Value = "'83.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").Value = "'2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("D27").Value = "'12.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'9.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("D31").Value = '4.006.46'
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").Value = "'7.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("D33").Value = "'2.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").Value = "'31.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("D35").Value = '3.832.81'
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("E36").Value = '  -2.24%  '
$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").Value = "'1.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = "'0.140"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").Value = "'5.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.57%  '
$ws.Range("D40").Value = "'3.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.79%  '
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = "'0.312"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D43").Value = "'427.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = "'47.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.60%  '
$ws.Range("D47").Value = "'8.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("D48").Value = "'27.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("D49").Value = "'143.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("E50").Value = '  +3.00%  '
$ws.Range("D51").Value = "'39.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.54%  '
